$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 0.188
$ws.Range("C2").Value = 0.5679999999999999
$ws.Range("J2").Value = 0.012
$ws.Range("O2").Value = 0.004
$ws.Range("P2").Value = 0.132
$ws.Range("S2").Value = 0.096
$ws.Range("B3").Value = 0.0136986301369863
$ws.Range("C3").Value = 0.02054794520547945
$ws.Range("J3").Value = 0.00684931506849315
$ws.Range("P3").Value = 0.7876712328767124
$ws.Range("S3").Value = 0.1712328767123288
$ws.Range("J4").Value = 0.02173913043478261
$ws.Range("P4").Value = 0.6956521739130435
$ws.Range("S4").Value = 0.2826086956521739
$ws.Range("B6").Value = 0.04081632653061224
$ws.Range("D6").Value = 0.00816326530612245
$ws.Range("F6").Value = 0.06122448979591837
$ws.Range("J6").Value = 0.273469387755102
$ws.Range("O6").Value = 0.004081632653061225
$ws.Range("Q6").Value = 0.1224489795918367
$ws.Range("R6").Value = 0.1061224489795918
$ws.Range("S6").Value = 0.3836734693877551
$ws.Range("B7").Value = 0.08571428571428572
$ws.Range("D7").Value = 0.02857142857142857
$ws.Range("F7").Value = 0.04761904761904762
$ws.Range("J7").Value = 0.09523809523809523
$ws.Range("O7").Value = 0.02380952380952381
$ws.Range("Q7").Value = 0.2
$ws.Range("R7").Value = 0.05714285714285714
$ws.Range("S7").Value = 0.4619047619047619
$ws.Range("B8").Value = 0.07874015748031496
$ws.Range("D8").Value = 0.007874015748031496
$ws.Range("F8").Value = 0.05511811023622047
$ws.Range("J8").Value = 0.1377952755905512
$ws.Range("O8").Value = 0.01181102362204724
$ws.Range("Q8").Value = 0.1830708661417323
$ws.Range("R8").Value = 0.09251968503937008
$ws.Range("S8").Value = 0.4330708661417323
$ws.Range("B9").Value = 0.09433962264150944
$ws.Range("F9").Value = 0.1069182389937107
$ws.Range("J9").Value = 0.1446540880503145
$ws.Range("O9").Value = 0.01257861635220126
$ws.Range("Q9").Value = 0.1320754716981132
$ws.Range("R9").Value = 0.07547169811320754
$ws.Range("S9").Value = 0.4339622641509434
$ws.Range("B10").Value = 0.0867579908675799
$ws.Range("D10").Value = 0.02663622526636225
$ws.Range("E10").Value = 0.0015220700152207
$ws.Range("F10").Value = 0.06773211567732115
$ws.Range("J10").Value = 0.134703196347032
$ws.Range("O10").Value = 0.0076103500761035
$ws.Range("Q10").Value = 0.2100456621004566
$ws.Range("R10").Value = 0.08447488584474885
$ws.Range("S10").Value = 0.380517503805175
$ws.Range("G11").Value = 0.1166180758017493
$ws.Range("J11").Value = 0.1049562682215743
$ws.Range("K11").Value = 0.1720116618075802
$ws.Range("L11").Value = 0.5772594752186589
$ws.Range("S11").Value = 0.02915451895043732
$ws.Range("G12").Value = 0.75
$ws.Range("J12").Value = 0.1666666666666667
$ws.Range("K12").Value = 0.004901960784313725
$ws.Range("L12").Value = 0.009803921568627451
$ws.Range("S12").Value = 0.06862745098039216
$ws.Range("G13").Value = 0.6666666666666666
$ws.Range("J13").Value = 0.3111111111111111
$ws.Range("S13").Value = 0.02222222222222222
$ws.Range("F15").Value = 0.02463054187192118
$ws.Range("H15").Value = 0.2216748768472906
$ws.Range("I15").Value = 0.06896551724137931
$ws.Range("J15").Value = 0.354679802955665
$ws.Range("K15").Value = 0.06403940886699508
$ws.Range("M15").Value = 0.01477832512315271
$ws.Range("O15").Value = 0.0541871921182266
$ws.Range("S15").Value = 0.1970443349753695
$ws.Range("F16").Value = 0.02298850574712644
$ws.Range("H16").Value = 0.1666666666666667
$ws.Range("I16").Value = 0.05747126436781609
$ws.Range("J16").Value = 0.396551724137931
$ws.Range("K16").Value = 0.09195402298850575
$ws.Range("M16").Value = 0.03448275862068965
$ws.Range("N16").Value = 0.005747126436781609
$ws.Range("O16").Value = 0.06896551724137931
$ws.Range("S16").Value = 0.1551724137931035
$ws.Range("F17").Value = 0.02838427947598253
$ws.Range("H17").Value = 0.1965065502183406
$ws.Range("I17").Value = 0.07205240174672489
$ws.Range("J17").Value = 0.4039301310043668
$ws.Range("K17").Value = 0.1222707423580786
$ws.Range("M17").Value = 0.01528384279475982
$ws.Range("O17").Value = 0.04803493449781659
$ws.Range("S17").Value = 0.1135371179039301
$ws.Range("F18").Value = 0.02898550724637681
$ws.Range("H18").Value = 0.1932367149758454
$ws.Range("I18").Value = 0.09178743961352658
$ws.Range("J18").Value = 0.3429951690821256
$ws.Range("K18").Value = 0.1449275362318841
$ws.Range("M18").Value = 0.01449275362318841
$ws.Range("O18").Value = 0.06280193236714976
$ws.Range("S18").Value = 0.1207729468599034
$ws.Range("F19").Value = 0.0186706497386109
$ws.Range("H19").Value = 0.2255414488424197
$ws.Range("I19").Value = 0.06422703510082151
$ws.Range("J19").Value = 0.3622106049290515
$ws.Range("K19").Value = 0.1239731142643764
$ws.Range("M19").Value = 0.01941747572815534
$ws.Range("O19").Value = 0.07617625093353249
$ws.Range("S19").Value = 0.1097834204630321
